$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reassign match data (columns F:V) among rows that were reordered within same-date groups ---
# (Columns A-E: Indice/pais/torneio/temporada/data_partida are unchanged for these rows)

# Row 29
$ws.Range("F29").Value = "Ilirija"
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = "Rudar"
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 1.81
$ws.Range("K29").Value = "05/08/2023 05:42"
$ws.Range("L29").Value = 2.13
$ws.Range("M29").Value = "30/08/2023 16:51"
$ws.Range("N29").Value = 3.51
$ws.Range("O29").Value = "05/08/2023 05:42"
$ws.Range("P29").Value = 3.62
$ws.Range("Q29").Value = "30/08/2023 16:21"
$ws.Range("R29").Value = 3.42
$ws.Range("S29").Value = "05/08/2023 05:42"
$ws.Range("T29").Value = 2.95
$ws.Range("U29").Value = "30/08/2023 16:51"
$ws.Range("V29").Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-rudar/pv9TMWW7/"

# Row 30
$ws.Range("F30").Value = "Fuzinar"
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = "Grosuplje"
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 3.24
$ws.Range("K30").Value = "04/08/2023 05:42"
$ws.Range("L30").Value = 4.09
$ws.Range("M30").Value = "30/08/2023 15:38"
$ws.Range("N30").Value = 3.15
$ws.Range("O30").Value = "04/08/2023 05:42"
$ws.Range("P30").Value = 3.69
$ws.Range("Q30").Value = "30/08/2023 15:38"
$ws.Range("R30").Value = 2.03
$ws.Range("S30").Value = "04/08/2023 05:42"
$ws.Range("T30").Value = 1.75
$ws.Range("U30").Value = "30/08/2023 15:38"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/slovenia/2-snl/fuzinar-grosuplje/OMD8pzPE/"

# Row 31
$ws.Range("F31").Value = "Triglav"
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = "Dravinja"
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = 2.41
$ws.Range("K31").Value = "03/08/2023 23:25"
$ws.Range("L31").Value = 1.91
$ws.Range("M31").Value = "30/08/2023 16:02"
$ws.Range("N31").Value = 3.34
$ws.Range("O31").Value = "03/08/2023 23:25"
$ws.Range("P31").Value = 3.65
$ws.Range("Q31").Value = "30/08/2023 16:02"
$ws.Range("R31").Value = 2.66
$ws.Range("S31").Value = "03/08/2023 23:25"
$ws.Range("T31").Value = 3.49
$ws.Range("U31").Value = "30/08/2023 16:02"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/slovenia/2-snl/triglav-dravinja/bN6HPYnk/"

# Row 62
$ws.Range("F62").Value = "Bilje"
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = "Triglav"
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = 2.38
$ws.Range("K62").Value = "20/09/2023 03:12"
$ws.Range("L62").Value = 2.69
$ws.Range("M62").Value = "21/09/2023 15:59"
$ws.Range("N62").Value = 3.24
$ws.Range("O62").Value = "20/09/2023 03:12"
$ws.Range("P62").Value = 3.06
$ws.Range("Q62").Value = "21/09/2023 15:59"
$ws.Range("R62").Value = 2.51
$ws.Range("S62").Value = "20/09/2023 03:12"
$ws.Range("T62").Value = 2.6
$ws.Range("U62").Value = "21/09/2023 15:59"
$ws.Range("V62").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-triglav/OKlbEQ09/"

# Row 63
$ws.Range("F63").Value = "Jadran Dekani"
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = "Beltinci"
$ws.Range("I63").Value = 3
$ws.Range("J63").Value = 2.74
$ws.Range("K63").Value = "20/09/2023 03:12"
$ws.Range("L63").Value = 3.03
$ws.Range("M63").Value = "21/09/2023 15:59"
$ws.Range("N63").Value = 3.07
$ws.Range("O63").Value = "20/09/2023 03:12"
$ws.Range("P63").Value = 3.18
$ws.Range("Q63").Value = "21/09/2023 15:59"
$ws.Range("R63").Value = 2.3
$ws.Range("S63").Value = "20/09/2023 03:12"
$ws.Range("T63").Value = 2.27
$ws.Range("U63").Value = "21/09/2023 14:40"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/slovenia/2-snl/jadran-dekani-beltinci/lSoABOpS/"

# Row 64
$ws.Range("F64").Value = "Nafta"
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = "Tolmin"
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1.32
$ws.Range("K64").Value = "20/09/2023 03:12"
$ws.Range("L64").Value = 1.47
$ws.Range("M64").Value = "21/09/2023 15:58"
$ws.Range("N64").Value = 4.8
$ws.Range("O64").Value = "20/09/2023 03:12"
$ws.Range("P64").Value = 4.17
$ws.Range("Q64").Value = "21/09/2023 15:59"
$ws.Range("R64").Value = 6.17
$ws.Range("S64").Value = "20/09/2023 03:12"
$ws.Range("T64").Value = 6.11
$ws.Range("U64").Value = "21/09/2023 15:59"
$ws.Range("V64").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nafta-tolmin/Aym6C4VL/"

# Row 73
$ws.Range("F73").Value = "ND Gorica"
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = "Fuzinar"
$ws.Range("I73").Value = 1
$ws.Range("J73").Value = 1.29
$ws.Range("K73").Value = "28/09/2023 02:42"
$ws.Range("L73").Value = 1.25
$ws.Range("M73").Value = "29/09/2023 13:45"
$ws.Range("N73").Value = 4.98
$ws.Range("O73").Value = "28/09/2023 02:42"
$ws.Range("P73").Value = 6.11
$ws.Range("Q73").Value = "29/09/2023 15:29"
$ws.Range("R73").Value = 6.88
$ws.Range("S73").Value = "28/09/2023 02:42"
$ws.Range("T73").Value = 7.81
$ws.Range("U73").Value = "29/09/2023 15:29"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nd-gorica-fuzinar/xQJeJCr6/"

# Row 74
$ws.Range("F74").Value = "Ilirija"
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = "Primorje"
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3.72
$ws.Range("K74").Value = "28/09/2023 02:42"
$ws.Range("L74").Value = 4.94
$ws.Range("M74").Value = "29/09/2023 15:20"
$ws.Range("N74").Value = 3.4
$ws.Range("O74").Value = "28/09/2023 02:42"
$ws.Range("P74").Value = 3.8
$ws.Range("Q74").Value = "29/09/2023 15:20"
$ws.Range("R74").Value = 1.79
$ws.Range("S74").Value = "28/09/2023 02:42"
$ws.Range("T74").Value = 1.61
$ws.Range("U74").Value = "29/09/2023 15:20"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-primorje/f1UjKhTa/"

# Row 76
$ws.Range("F76").Value = "Rudar"
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = "Triglav"
$ws.Range("I76").Value = 2
$ws.Range("J76").Value = 2.29
$ws.Range("K76").Value = "29/09/2023 02:42"
$ws.Range("L76").Value = 2.25
$ws.Range("M76").Value = "30/09/2023 15:29"
$ws.Range("N76").Value = 3.2
$ws.Range("O76").Value = "29/09/2023 02:42"
$ws.Range("P76").Value = 3.37
$ws.Range("Q76").Value = "30/09/2023 15:25"
$ws.Range("R76").Value = 2.64
$ws.Range("S76").Value = "29/09/2023 02:42"
$ws.Range("T76").Value = 2.9
$ws.Range("U76").Value = "30/09/2023 15:24"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/slovenia/2-snl/rudar-triglav/jepXGJM7/"

# Row 77
$ws.Range("F77").Value = "NK Bistrica"
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = "Tabor Sezana"
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1.41
$ws.Range("K77").Value = "29/09/2023 02:42"
$ws.Range("L77").Value = 1.34
$ws.Range("M77").Value = "30/09/2023 15:24"
$ws.Range("N77").Value = 4.29
$ws.Range("O77").Value = "29/09/2023 02:42"
$ws.Range("P77").Value = 5.2
$ws.Range("Q77").Value = "30/09/2023 15:29"
$ws.Range("R77").Value = 5.33
$ws.Range("S77").Value = "29/09/2023 02:42"
$ws.Range("T77").Value = 6.9
$ws.Range("U77").Value = "30/09/2023 15:29"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bistrica-tabor-sezana/tbQnLYDg/"

# Row 78
$ws.Range("F78").Value = "Jadran Dekani"
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = "Grosuplje"
$ws.Range("I78").Value = 1
$ws.Range("J78").Value = 3.27
$ws.Range("K78").Value = "29/09/2023 02:42"
$ws.Range("L78").Value = 4
$ws.Range("M78").Value = "30/09/2023 15:20"
$ws.Range("N78").Value = 3.13
$ws.Range("O78").Value = "29/09/2023 02:42"
$ws.Range("P78").Value = 3.26
$ws.Range("Q78").Value = "30/09/2023 15:20"
$ws.Range("R78").Value = 1.99
$ws.Range("S78").Value = "29/09/2023 02:42"
$ws.Range("T78").Value = 1.89
$ws.Range("U78").Value = "30/09/2023 15:20"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/slovenia/2-snl/jadran-dekani-grosuplje/rkXIvNEE/"

# Row 79
$ws.Range("F79").Value = "Nafta"
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = "Bilje"
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = 1.72
$ws.Range("K79").Value = "29/09/2023 02:42"
$ws.Range("L79").Value = 1.78
$ws.Range("M79").Value = "30/09/2023 15:23"
$ws.Range("N79").Value = 3.77
$ws.Range("O79").Value = "29/09/2023 02:42"
$ws.Range("P79").Value = 4.06
$ws.Range("Q79").Value = "30/09/2023 15:23"
$ws.Range("R79").Value = 3.65
$ws.Range("S79").Value = "29/09/2023 02:42"
$ws.Range("T79").Value = 3.61
$ws.Range("U79").Value = "30/09/2023 15:23"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nafta-bilje/4bYEu3a8/"

# Row 82
$ws.Range("F82").Value = "Grosuplje"
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = "Dravinja"
$ws.Range("I82").Value = 2
$ws.Range("J82").Value = 1.34
$ws.Range("K82").Value = "06/10/2023 02:42"
$ws.Range("L82").Value = 1.29
$ws.Range("M82").Value = "07/10/2023 15:03"
$ws.Range("N82").Value = 4.49
$ws.Range("O82").Value = "06/10/2023 02:42"
$ws.Range("P82").Value = 5.25
$ws.Range("Q82").Value = "07/10/2023 15:19"
$ws.Range("R82").Value = 6.24
$ws.Range("S82").Value = "06/10/2023 02:42"
$ws.Range("T82").Value = 8.34
$ws.Range("U82").Value = "07/10/2023 15:19"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-dravinja/f1jMZXSn/"

# Row 83
$ws.Range("F83").Value = "Beltinci"
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = "NK Bistrica"
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1.71
$ws.Range("K83").Value = "06/10/2023 02:42"
$ws.Range("L83").Value = 1.52
$ws.Range("M83").Value = "07/10/2023 10:57"
$ws.Range("N83").Value = 3.56
$ws.Range("O83").Value = "06/10/2023 02:42"
$ws.Range("P83").Value = 4.11
$ws.Range("Q83").Value = "07/10/2023 13:35"
$ws.Range("R83").Value = 3.79
$ws.Range("S83").Value = "06/10/2023 02:42"
$ws.Range("T83").Value = 5.27
$ws.Range("U83").Value = "07/10/2023 10:57"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/slovenia/2-snl/beltinci-bistrica/Qo0juf5P/"

# Row 84
$ws.Range("F84").Value = "Bilje"
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = "Jadran Dekani"
$ws.Range("I84").Value = 2
$ws.Range("J84").Value = 2.04
$ws.Range("K84").Value = "06/10/2023 02:42"
$ws.Range("L84").Value = 2.34
$ws.Range("M84").Value = "07/10/2023 15:15"
$ws.Range("N84").Value = 3.22
$ws.Range("O84").Value = "06/10/2023 02:42"
$ws.Range("P84").Value = 3.34
$ws.Range("Q84").Value = "07/10/2023 15:15"
$ws.Range("R84").Value = 3.06
$ws.Range("S84").Value = "06/10/2023 02:42"
$ws.Range("T84").Value = 2.8
$ws.Range("U84").Value = "07/10/2023 15:15"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-jadran-dekani/YikQYirg/"

# Row 85
$ws.Range("F85").Value = "Rudar"
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = "Nafta"
$ws.Range("I85").Value = 5
$ws.Range("J85").Value = 2.62
$ws.Range("K85").Value = "07/10/2023 02:42"
$ws.Range("L85").Value = 3.04
$ws.Range("M85").Value = "08/10/2023 15:01"
$ws.Range("N85").Value = 3.22
$ws.Range("O85").Value = "07/10/2023 02:42"
$ws.Range("P85").Value = 3.36
$ws.Range("Q85").Value = "08/10/2023 15:01"
$ws.Range("R85").Value = 2.29
$ws.Range("S85").Value = "07/10/2023 02:42"
$ws.Range("T85").Value = 2.18
$ws.Range("U85").Value = "08/10/2023 15:01"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/slovenia/2-snl/rudar-nafta/jTqVXBca/"

# Row 86
$ws.Range("F86").Value = "Tabor Sezana"
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = "Ilirija"
$ws.Range("I86").Value = 3
$ws.Range("J86").Value = 2.97
$ws.Range("K86").Value = "07/10/2023 02:42"
$ws.Range("L86").Value = 3.94
$ws.Range("M86").Value = "08/10/2023 13:47"
$ws.Range("N86").Value = 3.28
$ws.Range("O86").Value = "07/10/2023 02:42"
$ws.Range("P86").Value = 3.92
$ws.Range("Q86").Value = "08/10/2023 15:01"
$ws.Range("R86").Value = 2.05
$ws.Range("S86").Value = "07/10/2023 02:42"
$ws.Range("T86").Value = 1.72
$ws.Range("U86").Value = "08/10/2023 13:47"
$ws.Range("V86").Value = "https://www.betexplorer.com/football/slovenia/2-snl/tabor-sezana-ilirija/OtM7GASO/"

# Row 87
$ws.Range("F87").Value = "Tolmin"
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = "NK Krka"
$ws.Range("I87").Value = 3
$ws.Range("J87").Value = 3.75
$ws.Range("K87").Value = "07/10/2023 02:42"
$ws.Range("L87").Value = 2.86
$ws.Range("M87").Value = "08/10/2023 15:27"
$ws.Range("N87").Value = 3.48
$ws.Range("O87").Value = "07/10/2023 02:42"
$ws.Range("P87").Value = 3.82
$ws.Range("Q87").Value = "08/10/2023 15:28"
$ws.Range("R87").Value = 1.74
$ws.Range("S87").Value = "07/10/2023 02:42"
$ws.Range("T87").Value = 2.11
$ws.Range("U87").Value = "08/10/2023 15:27"
$ws.Range("V87").Value = "https://www.betexplorer.com/football/slovenia/2-snl/tolmin-nk-krka/0viIzEDt/"

# Row 95
$ws.Range("F95").Value = "Ilirija"
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = "Beltinci"
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 3.48
$ws.Range("K95").Value = "14/10/2023 02:12"
$ws.Range("L95").Value = 5.99
$ws.Range("M95").Value = "15/10/2023 14:43"
$ws.Range("N95").Value = 3.44
$ws.Range("O95").Value = "14/10/2023 02:12"
$ws.Range("P95").Value = 4.19
$ws.Range("Q95").Value = "15/10/2023 14:43"
$ws.Range("R95").Value = 1.81
$ws.Range("S95").Value = "14/10/2023 02:12"
$ws.Range("T95").Value = 1.47
$ws.Range("U95").Value = "15/10/2023 14:41"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-beltinci/lOBelXst/"

# Row 96
$ws.Range("F96").Value = "Nafta"
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = "Triglav"
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1.72
$ws.Range("K96").Value = "14/10/2023 02:12"
$ws.Range("L96").Value = 1.6
$ws.Range("M96").Value = "15/10/2023 14:41"
$ws.Range("N96").Value = 3.57
$ws.Range("O96").Value = "14/10/2023 02:12"
$ws.Range("P96").Value = 4.03
$ws.Range("Q96").Value = "15/10/2023 14:41"
$ws.Range("R96").Value = 3.73
$ws.Range("S96").Value = "14/10/2023 02:12"
$ws.Range("T96").Value = 4.77
$ws.Range("U96").Value = "15/10/2023 14:41"
$ws.Range("V96").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nafta-triglav/UJrZWVC5/"

# Row 106
$ws.Range("F106").Value = "Grosuplje"
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = "Primorje"
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 2.62
$ws.Range("K106").Value = "08/08/2023 04:42"
$ws.Range("L106").Value = 2.54
$ws.Range("M106").Value = "24/10/2023 14:54"
$ws.Range("N106").Value = 3.06
$ws.Range("O106").Value = "08/08/2023 04:42"
$ws.Range("P106").Value = 3.05
$ws.Range("Q106").Value = "24/10/2023 14:59"
$ws.Range("R106").Value = 2.44
$ws.Range("S106").Value = "08/08/2023 04:42"
$ws.Range("T106").Value = 2.75
$ws.Range("U106").Value = "24/10/2023 14:54"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/slovenia/2-snl/grosuplje-primorje/pfcixCWf/"

# Row 107
$ws.Range("F107").Value = "Bilje"
$ws.Range("G107").Value = 5
$ws.Range("H107").Value = "Fuzinar"
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 1.89
$ws.Range("K107").Value = "23/10/2023 02:12"
$ws.Range("L107").Value = 2
$ws.Range("M107").Value = "24/10/2023 14:56"
$ws.Range("N107").Value = 3.58
$ws.Range("O107").Value = "23/10/2023 02:12"
$ws.Range("P107").Value = 3.68
$ws.Range("Q107").Value = "24/10/2023 14:58"
$ws.Range("R107").Value = 3.12
$ws.Range("S107").Value = "23/10/2023 02:12"
$ws.Range("T107").Value = 3.19
$ws.Range("U107").Value = "24/10/2023 14:58"
$ws.Range("V107").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-fuzinar/OpdeyWo1/"

# Row 108
$ws.Range("F108").Value = "Beltinci"
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = "Triglav"
$ws.Range("I108").Value = 2
$ws.Range("J108").Value = 2.13
$ws.Range("K108").Value = "08/08/2023 04:42"
$ws.Range("L108").Value = 1.49
$ws.Range("M108").Value = "25/10/2023 14:59"
$ws.Range("N108").Value = 3.25
$ws.Range("O108").Value = "08/08/2023 04:42"
$ws.Range("P108").Value = 4.51
$ws.Range("Q108").Value = "25/10/2023 14:59"
$ws.Range("R108").Value = 2.86
$ws.Range("S108").Value = "08/08/2023 04:42"
$ws.Range("T108").Value = 5.21
$ws.Range("U108").Value = "25/10/2023 14:59"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/slovenia/2-snl/beltinci-triglav/xUgqvY1r/"

# Row 109
$ws.Range("F109").Value = "Dravinja"
$ws.Range("G109").Value = 3
$ws.Range("H109").Value = "NK Krka"
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 3.2
$ws.Range("K109").Value = "24/10/2023 02:12"
$ws.Range("L109").Value = 2.88
$ws.Range("M109").Value = "25/10/2023 14:57"
$ws.Range("N109").Value = 3.36
$ws.Range("O109").Value = "24/10/2023 02:12"
$ws.Range("P109").Value = 3.39
$ws.Range("Q109").Value = "25/10/2023 14:57"
$ws.Range("R109").Value = 1.93
$ws.Range("S109").Value = "24/10/2023 02:12"
$ws.Range("T109").Value = 2.26
$ws.Range("U109").Value = "25/10/2023 14:57"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/slovenia/2-snl/dravinja-nk-krka/AwxjcD8K/"

# Row 110
$ws.Range("F110").Value = "Jadran Dekani"
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = "NK Bistrica"
$ws.Range("I110").Value = 4
$ws.Range("J110").Value = 2.31
$ws.Range("K110").Value = "24/10/2023 02:12"
$ws.Range("L110").Value = 2.05
$ws.Range("M110").Value = "25/10/2023 14:57"
$ws.Range("N110").Value = 3.11
$ws.Range("O110").Value = "24/10/2023 02:12"
$ws.Range("P110").Value = 3.32
$ws.Range("Q110").Value = "25/10/2023 14:58"
$ws.Range("R110").Value = 2.68
$ws.Range("S110").Value = "24/10/2023 02:12"
$ws.Range("T110").Value = 3.37
$ws.Range("U110").Value = "25/10/2023 14:58"
$ws.Range("V110").Value = "https://www.betexplorer.com/football/slovenia/2-snl/jadran-dekani-bistrica/WnwnbggE/"

# Row 111
$ws.Range("F111").Value = "Tolmin"
$ws.Range("G111").Value = 3
$ws.Range("H111").Value = "Tabor Sezana"
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 1.69
$ws.Range("K111").Value = "24/10/2023 02:12"
$ws.Range("L111").Value = 1.96
$ws.Range("M111").Value = "25/10/2023 14:50"
$ws.Range("N111").Value = 3.65
$ws.Range("O111").Value = "24/10/2023 02:12"
$ws.Range("P111").Value = 3.8
$ws.Range("Q111").Value = "25/10/2023 14:50"
$ws.Range("R111").Value = 3.78
$ws.Range("S111").Value = "24/10/2023 02:12"
$ws.Range("T111").Value = 3.2
$ws.Range("U111").Value = "25/10/2023 14:50"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/slovenia/2-snl/tolmin-tabor-sezana/GKhmwhHl/"

# Row 113
$ws.Range("F113").Value = "Primorje"
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = "Tabor Sezana"
$ws.Range("I113").Value = 1
$ws.Range("J113").Value = 1.19
$ws.Range("K113").Value = "27/10/2023 02:13"
$ws.Range("L113").Value = 1.28
$ws.Range("M113").Value = "28/10/2023 13:41"
$ws.Range("N113").Value = 5.89
$ws.Range("O113").Value = "27/10/2023 02:13"
$ws.Range("P113").Value = 5.51
$ws.Range("Q113").Value = "28/10/2023 13:42"
$ws.Range("R113").Value = 8.67
$ws.Range("S113").Value = "27/10/2023 02:13"
$ws.Range("T113").Value = 8.279999999999999
$ws.Range("U113").Value = "28/10/2023 13:42"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/slovenia/2-snl/primorje-tabor-sezana/UHwRkpJp/"

# Row 115
$ws.Range("F115").Value = "NK Bistrica"
$ws.Range("G115").Value = 4
$ws.Range("H115").Value = "Bilje"
$ws.Range("I115").Value = 3
$ws.Range("J115").Value = 1.92
$ws.Range("K115").Value = "27/10/2023 02:13"
$ws.Range("L115").Value = 2.13
$ws.Range("M115").Value = "28/10/2023 13:53"
$ws.Range("N115").Value = 3.44
$ws.Range("O115").Value = "27/10/2023 02:13"
$ws.Range("P115").Value = 3.39
$ws.Range("Q115").Value = "28/10/2023 13:53"
$ws.Range("R115").Value = 3.16
$ws.Range("S115").Value = "27/10/2023 02:13"
$ws.Range("T115").Value = 3.13
$ws.Range("U115").Value = "28/10/2023 13:53"
$ws.Range("V115").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bistrica-bilje/vXWgcSQA/"

# Row 118
$ws.Range("F118").Value = "ND Gorica"
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = "Tolmin"
$ws.Range("I118").Value = 3
$ws.Range("J118").Value = 1.27
$ws.Range("K118").Value = "28/10/2023 03:12"
$ws.Range("L118").Value = 1.32
$ws.Range("M118").Value = "29/10/2023 11:36"
$ws.Range("N118").Value = 5
$ws.Range("O118").Value = "28/10/2023 03:12"
$ws.Range("P118").Value = 4.94
$ws.Range("Q118").Value = "29/10/2023 12:03"
$ws.Range("R118").Value = 7.14
$ws.Range("S118").Value = "28/10/2023 03:12"
$ws.Range("T118").Value = 7.46
$ws.Range("U118").Value = "29/10/2023 11:36"
$ws.Range("V118").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nd-gorica-tolmin/ncz2e6eN/"

# Row 119
$ws.Range("F119").Value = "Ilirija"
$ws.Range("G119").Value = 4
$ws.Range("H119").Value = "Grosuplje"
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 3.17
$ws.Range("K119").Value = "28/10/2023 03:12"
$ws.Range("L119").Value = 3.45
$ws.Range("M119").Value = "29/10/2023 13:45"
$ws.Range("N119").Value = 3.25
$ws.Range("O119").Value = "28/10/2023 03:12"
$ws.Range("P119").Value = 3.4
$ws.Range("Q119").Value = "29/10/2023 13:45"
$ws.Range("R119").Value = 1.98
$ws.Range("S119").Value = "28/10/2023 03:12"
$ws.Range("T119").Value = 2
$ws.Range("U119").Value = "29/10/2023 13:45"
$ws.Range("V119").Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-grosuplje/hzybdntH/"

# --- Append 2 new rows (127, 128) with matching style (copy format from row 126) ---
$ws.Range("A126:V126").Copy() | Out-Null
$ws.Range("A127:V128").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 127
$ws.Range("A127").Value = 126
$ws.Range("B127").Value = "slovenia"
$ws.Range("C127").Value = "2-snl"
$ws.Range("D127").Value = "2023-2024"
$ws.Range("E127").Value = 45236.58333333334
$ws.Range("F127").Value = "NK Bistrica"
$ws.Range("G127").Value = 5
$ws.Range("H127").Value = "Rudar"
$ws.Range("I127").Value = 1
$ws.Range("J127").Value = 1.95
$ws.Range("K127").Value = "05/11/2023 02:12"
$ws.Range("L127").Value = 2.17
$ws.Range("M127").Value = "06/11/2023 13:56"
$ws.Range("N127").Value = 3.4
$ws.Range("O127").Value = "05/11/2023 02:12"
$ws.Range("P127").Value = 3.44
$ws.Range("Q127").Value = "06/11/2023 13:56"
$ws.Range("R127").Value = 3.1
$ws.Range("S127").Value = "05/11/2023 02:12"
$ws.Range("T127").Value = 3
$ws.Range("U127").Value = "06/11/2023 13:56"
$ws.Range("V127").Value = "https://www.betexplorer.com/football/slovenia/2-snl/bistrica-rudar/Q3NT95AG/"

# Row 128
$ws.Range("A128").Value = 127
$ws.Range("B128").Value = "slovenia"
$ws.Range("C128").Value = "2-snl"
$ws.Range("D128").Value = "2023-2024"
$ws.Range("E128").Value = 45236.58333333334
$ws.Range("F128").Value = "Ilirija"
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = "Bilje"
$ws.Range("I128").Value = 3
$ws.Range("J128").Value = 2.03
$ws.Range("K128").Value = "05/11/2023 02:12"
$ws.Range("L128").Value = 2.04
$ws.Range("M128").Value = "06/11/2023 13:57"
$ws.Range("N128").Value = 3.41
$ws.Range("O128").Value = "05/11/2023 02:12"
$ws.Range("P128").Value = 3.66
$ws.Range("Q128").Value = "06/11/2023 13:57"
$ws.Range("R128").Value = 2.92
$ws.Range("S128").Value = "05/11/2023 02:12"
$ws.Range("T128").Value = 3.12
$ws.Range("U128").Value = "06/11/2023 13:57"
$ws.Range("V128").Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-bilje/dUIPAoeA/"
